$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 38; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value()
    $serial = $current.ToOADate()
    if ($serial -eq 46060) {
        $cell.Value = 46061
    }
}
